# Update the "取得日時" (acquired timestamp) column for the listed rows
# from 2025-10-15 06:26:22 to 2025-10-15 06:34:35, matching the commit
# message "Append: 2025-10-15 06:34 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-15 06:34:35"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
